$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (title reflects date through 12-22 instead of 12-21)
$ws.Name = "Through 2022-12-22"

# Update the label for the December row
$ws.Range("A13").Value = "December (through 12-22)"

# Update December row (row 13) values
$ws.Range("B13").Value = 30
$ws.Range("C13").Value = 69
$ws.Range("D13").Value = 88
$ws.Range("E13").Value = 50
$ws.Range("F13").Value = 44
$ws.Range("G13").Value = 104
$ws.Range("H13").Value = 148
$ws.Range("I13").Value = 98

# Update Total row (row 14) values
$ws.Range("B14").Value = 321
$ws.Range("C14").Value = 632
$ws.Range("D14").Value = 909
$ws.Range("E14").Value = 732
$ws.Range("F14").Value = 578
$ws.Range("G14").Value = 1368
$ws.Range("H14").Value = 1791
$ws.Range("I14").Value = 1615
